# "Update countries & provincias Spain"
#
# Daily refresh of the COVID "paises" workbook:
#   1) Bump the "last updated" timestamp (A1).
#   2) Update case counters (Casos totales / Nuevos casos / Casos activos /
#      Recuperados / Casos criticos / Muertes hoy / Muertes) for the
#      countries whose figures moved.
#   3) A few neighbouring countries swapped sort order (their totals in
#      column B are tied/very close), so their row data - and in two cases
#      the country name itself - are exchanged between the two rows to
#      keep the sheet's descending-by-total-cases ordering intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Timestamp ---------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Junio de 2020 a las 11:24"

# --- 2) Plain data refreshes (no reordering) ------------------------------
# Row 25
$ws.Range("B25").Value = 59918
$ws.Range("C25").Value = 99
$ws.Range("D25").Value = 16547
$ws.Range("E25").Value = 33721
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 9650

# Row 34
$ws.Range("B34").Value = 37420
$ws.Range("C34").Value = 1014
$ws.Range("D34").Value = 13776
$ws.Range("E34").Value = 21553
$ws.Range("G34").Value = 43
$ws.Range("H34").Value = 2091

# Row 52
$ws.Range("B52").Value = 17078
$ws.Range("C52").Value = 14
$ws.Range("D52").Value = 16012
$ws.Range("E52").Value = 389
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 677

# Row 61
$ws.Range("D61").Value = 6421
$ws.Range("E61").Value = 4277
$ws.Range("G61").Value = 10
$ws.Range("H61").Value = 395

# Row 69
$ws.Range("B69").Value = 8445
$ws.Range("C69").Value = 43
$ws.Range("D69").Value = 7311
$ws.Range("E69").Value = 1014
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 120

# Row 191 (Namibia)
$ws.Range("B191").Value = 32
$ws.Range("C191").Value = 1
$ws.Range("E191").Value = 15

# --- 3) Uganda / San Marino swap (rows 137-138) ---------------------------
$ws.Range("A137").Value = "Uganda"
$ws.Range("C137").Value = 8
$ws.Range("D137").Value = 199
$ws.Range("E137").Value = 495
$ws.Range("H137").Value = 0

$ws.Range("A138").Value = "San Marino"
$ws.Range("B138").Value = 694
$ws.Range("D138").Value = 520
$ws.Range("E138").Value = 132
$ws.Range("H138").Value = 42

# --- 4) Small-island nations reshuffle (rows 206-214) ---------------------
# Groenlandia <-> Islas Malvinas
$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"

# Santa Sede <-> Islas Turcas y Caicos
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

# Montserrat <-> Seychelles
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Islas Virgenes Britanicas <-> Papua Nueva Guinea
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
